# The edit re-orders the two observation records currently stored in rows
# 18 and 19: the "Garnlav" (lichen) record that used to be row 18 becomes
# row 19, and the "Tjäder" (bird) record that used to be row 19 becomes
# row 18. Column A..AY hold the rest of each record's fields, several of
# which (I, J, K, L, M, N, AC, AF) are only populated for one of the two
# species, so cells that must become blank are cleared and cells that must
# newly appear (blank placeholders) are created by copying an already
# blank cell onto them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: materialise the new blank placeholder cells ---------------
# (I18 is blank in both the source and target state, so it is a safe,
# untouched donor for creating brand-new empty cells elsewhere.)
$ws.Range("I18").Copy($ws.Range("L18"))
$ws.Range("I18").Copy($ws.Range("J19"))
$ws.Range("I18").Copy($ws.Range("AF19"))

# --- Step 2: drop the cells that must no longer exist -------------------
$ws.Range("J18").ClearContents()
$ws.Range("AC18").ClearContents()
$ws.Range("AF18").ClearContents()
$ws.Range("L19").ClearContents()
$ws.Range("M19").ClearContents()

# --- Step 3: write the swapped record data ------------------------------
# New row 18 = old row 19 data ("Tjäder")
$ws.Range("A18").Value = 131187791
$ws.Range("B18").Value = 57073
$ws.Range("D18").Value = "LC"
$ws.Range("E18").Value = 100138
$ws.Range("F18").Value = "Tjäder"
$ws.Range("G18").Value = "Tetrao urogallus"
$ws.Range("H18").Value = "Linnaeus, 1758"
$ws.Range("M18").Value = "färsk spillning"
$ws.Range("P18").Value = "Svatå, Dlr"
$ws.Range("Q18").Value = 511301
$ws.Range("R18").Value = 6697864

# New row 19 = old row 18 data ("Garnlav")
$ws.Range("A19").Value = 131187762
$ws.Range("B19").Value = 79244
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6425
$ws.Range("F19").Value = "Garnlav"
$ws.Range("G19").Value = "Alectoria sarmentosa"
$ws.Range("H19").Value = "(Ach.) Ach."
$ws.Range("P19").Value = "Svartå, Dlr"
$ws.Range("Q19").Value = 511511
$ws.Range("R19").Value = 6697866
$ws.Range("AC19").Value = "På äldre tall."
